$d = $word.ActiveDocument

# Helper: append a new run of text to the end of the given paragraph,
# applying the same Times New Roman / en-US formatting used throughout
# the document.
function Append-FormattedRun {
    param($range, [string]$text)
    $range.Collapse(0)
    $range.InsertAfter($text)
    $range.Font.Name = "Times New Roman"
    $range.Font.NameBi = "Times New Roman"
    $range.LanguageID = "en-US"
}

# The last paragraph currently ends with "...rethink our whole approach to
# creating the program." Append four additional runs of text to that same
# paragraph.
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$r = $lastPara.Range

Append-FormattedRun $r " We had to put into consideration the difference in"
Append-FormattedRun $r " data"
Append-FormattedRun $r " structures"
Append-FormattedRun $r " when reimplementing the code."

# Add a brand-new paragraph after the current last paragraph, then fill it
# with the formatted text.
$r.Collapse(0)
$r.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$r = $newPara.Range
Append-FormattedRun $r "We realized the most effective want to store the data was to use vectors and a HashMap. The HashMap was implemented inside the route.cpp file where we were able to read the csv files provided and extract the appropriate columns."

# Add a second brand-new paragraph after that one.
$r.Collapse(0)
$r.InsertParagraphAfter()
$newPara2 = $d.Paragraphs.Item($d.Paragraphs.Count)
$r = $newPara2.Range
Append-FormattedRun $r "In conclusion, our program ended with taking an input file called test.txt as a parameter and it would provide the relevant information of the start and destination location in the test file. The relevant information will then be written to an output file which has the names of the start and end locations, which was also a text file."
